# The sheet gained one new weekly price record. It is inserted as the new
# row 27 (pushing the former rows 27-100 down to 28-101), so the used range
# grows from A1:R100 to A1:R101.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(27).Insert()

$ws.Range("A27").Value = 6
$ws.Range("B27").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C27").Value = "Metropolitana"
$ws.Range("D27").Value = 44414
$ws.Range("E27").Value = 13
$ws.Range("F27").Value = 100112022
$ws.Range("G27").Value = "Arveja Verde"
$ws.Range("H27").Value = "Perfection"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 200
$ws.Range("K27").Value = 35000
$ws.Range("L27").Value = 36000
$ws.Range("M27").Value = 35600
$ws.Range("N27").Value = '$/malla 25 kilos'
$ws.Range("O27").Value = "Provincia de Huasco"
$ws.Range("P27").Value = 1424
$ws.Range("Q27").Value = 25
$ws.Range("R27").Value = "Hortaliza"
